$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19: new data appended after the last existing row (18)
# Force the date-like string to stay as literal text (matching the
# other text-dates already in column A) instead of being auto-converted
# to a date serial number, then reset the style back to the default
# so no stray style index is left on the cell.
$ws.Cells.Item(19, 1).NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = "05/03/2025"
$ws.Cells.Item(19, 1).Style = "Normal"

$ws.Cells.Item(19, 2).Value = 512.6880000000019
$ws.Cells.Item(19, 3).Value = 0.09752520051181189
$ws.Cells.Item(19, 4).Value = 50
